# Weekly update: a new price record (week of 44509) is inserted at the top
# of the data block (row 199), pushing the existing rows 199-205 down to
# 200-206.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 199; this shifts rows 199:205
# down to 200:206 and keeps their contents/formatting intact.
$ws.Rows.Item(199).Insert()

# Populate the newly inserted row 199 with the new weekly record.
$ws.Range('A199').Value = 8
$ws.Range('B199').Value = 'Terminal La Palmera de La Serena'
$ws.Range('C199').Value = 'Coquimbo'
$ws.Range('D199').Value = 44509
$ws.Range('E199').Value = 4
$ws.Range('F199').Value = 100112032
$ws.Range('G199').Value = 'Zapallo italiano'
$ws.Range('H199').Value = 'Sin especificar'
$ws.Range('I199').Value = 'Primera'
$ws.Range('J199').Value = 560
$ws.Range('K199').Value = 10000
$ws.Range('L199').Value = 11000
$ws.Range('M199').Value = 10500
$ws.Range('N199').Value = '$/caja 70 unidades'
$ws.Range('O199').Value = 'Provincia de Limarí'
$ws.Range('P199').Value = 150
$ws.Range('Q199').Value = 70
$ws.Range('R199').Value = 'Hortaliza'
